# Generate Report for Handoff
# Applies the localization-status refresh: new source GUID, new xliff hash,
# updated timestamps, cleared "Latest Target File"/"Latest Handback File"
# columns (and their hyperlink on the old source filename), plus the
# matching column-width adjustments.

$wb = $excel.ActiveWorkbook

$oldGuid = "fdfe4cc4-e100-4b84-b26c-6f430fac2abf"
$newGuid = "120cec05-8e5f-460e-a85a-33d820ad377a"
$newHash = "34436264f6623362b94f428279ea481df536ad00"

# --- Overview sheet --------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "$newGuid.md"
$ov.Range("B2").Value = "e2e\$newGuid.md"
$ov.Range("G2").Value = "2016-08-30 04:59:33"

foreach ($hl in $ov.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

$ov.Columns.Item(1).ColumnWidth = 39.17

# --- zh-cn / de-de sheets ---------------------------------------------
$langs = @{ "zh-cn" = "2016-08-30 04:59:29"; "de-de" = "2016-08-30 04:59:33" }

foreach ($langName in $langs.Keys) {
    $ws = $wb.Worksheets.Item($langName)
    $handoffDate = $langs[$langName]

    $ws.Range("A2").Value = "$newGuid.md"
    $ws.Range("G2").Value = "$newGuid.$newHash.$langName.xlf"
    $ws.Range("H2").Value = $handoffDate
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # Update the surviving hyperlink's display text (A2 keeps its link;
    # the old one sitting on I2 -- the stale "Latest Target File" -- is
    # removed entirely below). The cell text and the hyperlink's display
    # text are independent properties in this model, so both are set.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$2') {
            $hl.TextToDisplay = "$newGuid.md"
        }
    }

    for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
        $hl = $ws.Hyperlinks.Item($i)
        if ($hl.Range.Address() -eq '$I$2') {
            $hl.Delete()
        }
    }

    # I2 loses its hyperlink formatting along with the hyperlink itself.
    $ws.Range("I2").Style = "Normal"

    $ws.Columns.Item(1).ColumnWidth = 39.17
    $ws.Columns.Item(9).ColumnWidth = 17.8
    $ws.Columns.Item(10).ColumnWidth = 20.8
}
